# GBDS DECEMBER FILES 2025 - fliqlo@GBDS
# Rename the "PE, NOVEMBER" sheet to "PE, DECEMBER" and fill in the first
# purchase-entry row (row 7) with December's data.

$wb = $excel.ActiveWorkbook

# --- Rename the worksheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("PE, NOVEMBER")
$ws.Name = "PE, DECEMBER"

# --- Update the Print_Area defined name to point at the renamed sheet ------
$printArea = $wb.Names.Item("PE, DECEMBER!Print_Area")
$printArea.RefersTo = "='PE, DECEMBER'!#REF!"

# --- Restore the active selection shown in the sheet view ------------------
$ws.Range("I8").Select()

# --- Fill in row 7 with the December purchase entry -------------------------
$ws.Range("C7").Value = 45973
$ws.Range("G7").Value = 518184389
$ws.Range("I7").Formula = "=1353132-110674.04"
